{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The document is a multiplication-drill worksheet: a centered date\n// paragraph followed by a 20-row x 5-column table whose rows 1, 5, 10,\n// 15 and 20 (1-indexed) hold \"three-digit x one-digit\" problems such as\n// \"251\u00d77=1757\"; the other rows are blank answer rows. This script\n// updates the date and all 25 problem cells to their new values using\n// exact text search-and-replace, which is robust to the table/paragraph\n// structure staying otherwise untouched.\n\nconst replacements = [\n  [\"2025-11-07 Friday\", \"2025-11-08 Saturday\"],\n  [\"251\u00d77=1757\", \"434\u00d79=3906\"],\n  [\"532\u00d77=3724\", \"714\u00d75=3570\"],\n  [\"966\u00d73=2898\", \"439\u00d73=1317\"],\n  [\"768\u00d79=6912\", \"449\u00d77=3143\"],\n  [\"692\u00d77=4844\", \"673\u00d76=4038\"],\n  [\"398\u00d79=3582\", \"435\u00d72=870\"],\n  [\"554\u00d74=2216\", \"603\u00d74=2412\"],\n  [\"514\u00d76=3084\", \"385\u00d74=1540\"],\n  [\"868\u00d73=2604\", \"915\u00d77=6405\"],\n  [\"820\u00d75=4100\", \"495\u00d74=1980\"],\n  [\"553\u00d74=2212\", \"263\u00d74=1052\"],\n  [\"412\u00d73=1236\", \"510\u00d75=2550\"],\n  [\"215\u00d78=1720\", \"758\u00d74=3032\"],\n  [\"529\u00d77=3703\", \"267\u00d74=1068\"],\n  [\"901\u00d73=2703\", \"129\u00d74=516\"],\n  [\"298\u00d75=1490\", \"266\u00d74=1064\"],\n  [\"963\u00d75=4815\", \"572\u00d72=1144\"],\n  [\"370\u00d77=2590\", \"179\u00d72=358\"],\n  [\"836\u00d74=3344\", \"434\u00d75=2170\"],\n  [\"437\u00d78=3496\", \"765\u00d78=6120\"],\n  [\"914\u00d79=8226\", \"197\u00d78=1576\"],\n  [\"166\u00d75=830\", \"835\u00d75=4175\"],\n  [\"184\u00d78=1472\", \"195\u00d74=780\"],\n  [\"228\u00d75=1140\", \"553\u00d72=1106\"],\n  [\"892\u00d79=8028\", \"479\u00d74=1916\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found, cannot apply replacement: \"${oldText}\"`);\n  }\n\n  // Every value in this worksheet is unique, but replace all hits just\n  // in case (mirrors \"Replace All\" semantics).\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# The document is open as $word.ActiveDocument ($d below).\n#\n# Replaces the worksheet date and the 25 three-digit x one-digit\n# multiplication problems with their updated values, matching the\n# target OOXML diff exactly.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-11-07 Friday', '2025-11-08 Saturday'),\n    @('251\u00d77=1757', '434\u00d79=3906'),\n    @('532\u00d77=3724', '714\u00d75=3570'),\n    @('966\u00d73=2898', '439\u00d73=1317'),\n    @('768\u00d79=6912', '449\u00d77=3143'),\n    @('692\u00d77=4844', '673\u00d76=4038'),\n    @('398\u00d79=3582', '435\u00d72=870'),\n    @('554\u00d74=2216', '603\u00d74=2412'),\n    @('514\u00d76=3084', '385\u00d74=1540'),\n    @('868\u00d73=2604', '915\u00d77=6405'),\n    @('820\u00d75=4100', '495\u00d74=1980'),\n    @('553\u00d74=2212', '263\u00d74=1052'),\n    @('412\u00d73=1236', '510\u00d75=2550'),\n    @('215\u00d78=1720', '758\u00d74=3032'),\n    @('529\u00d77=3703', '267\u00d74=1068'),\n    @('901\u00d73=2703', '129\u00d74=516'),\n    @('298\u00d75=1490', '266\u00d74=1064'),\n    @('963\u00d75=4815', '572\u00d72=1144'),\n    @('370\u00d77=2590', '179\u00d72=358'),\n    @('836\u00d74=3344', '434\u00d75=2170'),\n    @('437\u00d78=3496', '765\u00d78=6120'),\n    @('914\u00d79=8226', '197\u00d78=1576'),\n    @('166\u00d75=830', '835\u00d75=4175'),\n    @('184\u00d78=1472', '195\u00d74=780'),\n    @('228\u00d75=1140', '553\u00d72=1106'),\n    @('892\u00d79=8028', '479\u00d74=1916')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $ok = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Text not found, cannot apply replacement: '$oldText'\"\n    }\n}\n\n"}
